$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("E2").Value = 0.2669486146347392
$ws.Range("F2").Value = 0.3223287978607803
$ws.Range("G2").Value = 0.1145363079130303
$ws.Range("I2").Value = 0.03409398221989002
$ws.Range("M2").Value = 0.002190037553847609
$ws.Range("O2").Value = 0.001730587257784922
$ws.Range("P2").Value = 0.1549376786564502
$ws.Range("R2").Value = 0.05063323737114604
$ws.Range("T2").Value = 0.004712950693591489
$ws.Range("U2").Value = 0.01793159680657433
$ws.Range("X2").Value = 0.004961187174855878
$ws.Range("Z2").Value = 0.0209376856517416
$ws.Range("AD2").Value = 0.004057336205568327
$ws.Range("D3").Value = 0.007761633782800323
$ws.Range("E3").Value = 0.1742607829468068
$ws.Range("F3").Value = 0.3222927622629421
$ws.Range("G3").Value = 0.08400608921634399
$ws.Range("H3").Value = 0.004448714216209044
$ws.Range("I3").Value = 0.01652512429871952
$ws.Range("K3").Value = 0.01834455015247031
$ws.Range("L3").Value = 0.0005118182862332179
$ws.Range("M3").Value = 0.01749377673842456
$ws.Range("O3").Value = 0.01832592757464277
$ws.Range("P3").Value = 0.1677479890775728
$ws.Range("R3").Value = 0.07001026517306112
$ws.Range("S3").Value = 0.01236221865580689
$ws.Range("U3").Value = 0.04726059435227542
$ws.Range("X3").Value = 0.009741762482689313
$ws.Range("Z3").Value = 0.02036266458658111
$ws.Range("AD3").Value = 0.007572888956261446
$ws.Range("AE3").Value = 0.0009704372401592813
$ws.Range("F4").Value = 0.1394334263860402
$ws.Range("H4").Value = 0.3065401045915333
$ws.Range("I4").Value = 0.0008241197903705752
$ws.Range("J4").Value = 0.0003270086172496745
$ws.Range("K4").Value = 0.02318287153715201
$ws.Range("M4").Value = 0.008221677807233298
$ws.Range("O4").Value = 0.04151752705930281
$ws.Range("Q4").Value = 0.1189595447923376
$ws.Range("R4").Value = 0.1365039717685407
$ws.Range("S4").Value = 0.03537510659893097
$ws.Range("T4").Value = 0.01614491902403837
$ws.Range("U4").Value = 0.06822966446542007
$ws.Range("W4").Value = 0.03158269026507425
$ws.Range("Z4").Value = 0.02622547147363261
$ws.Range("AE4").Value = 0.01478435156407626
$ws.Range("AF4").Value = 0.03214754425906741
$ws.Range("F5").Value = 0.203605394793228
$ws.Range("H5").Value = 0.3082576828355844
$ws.Range("I5").Value = 0.01030339782758265
$ws.Range("J5").Value = 0.03773032822434538
$ws.Range("K5").Value = 0.02956120441574965
$ws.Range("M5").Value = 0.001035815403185701
$ws.Range("N5").Value = 0.005206590115208655
$ws.Range("O5").Value = 0.01029068963222245
$ws.Range("Q5").Value = 0.08460382708303539
$ws.Range("R5").Value = 0.08110055752858289
$ws.Range("S5").Value = 0.04432186595195839
$ws.Range("T5").Value = 0.005735244772170803
$ws.Range("U5").Value = 0.09827662833627222
$ws.Range("W5").Value = 0.009187524794622979
$ws.Range("X5").Value = 0.01746112712614082
$ws.Range("Z5").Value = 0.01199833877744786
$ws.Range("AA5").Value = 0.005634370942115148
$ws.Range("AB5").Value = 0.01168144614634196
$ws.Range("AF5").Value = 0.02400796529420475
$ws.Range("E6").Value = 0.2744254467974087
$ws.Range("F6").Value = 0.03806168675814205
$ws.Range("G6").Value = 0.1225653737173593
$ws.Range("I6").Value = 0.0341113696500616
$ws.Range("J6").Value = 0.02927396099844227
$ws.Range("K6").Value = 0.01891613647809891
$ws.Range("N6").Value = 0.01106019313681186
$ws.Range("P6").Value = 0.1383570672821516
$ws.Range("Q6").Value = 0.03391457080232953
$ws.Range("R6").Value = 0.1273635829965314
$ws.Range("T6").Value = 0.1378800257935012
$ws.Range("U6").Value = 0.0232594427984623
$ws.Range("V6").Value = 0.006495218258618761
$ws.Range("W6").Value = 0.0006755475637552843
$ws.Range("AE6").Value = 0.003640376968325262

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("E2").Value = 0.2669486146347392
$ws.Range("F2").Value = 0.5892774124955195
$ws.Range("G2").Value = 0.7038137204085497
$ws.Range("H2").Value = 0.7038137204085497
$ws.Range("I2").Value = 0.7379077026284397
$ws.Range("J2").Value = 0.7379077026284397
$ws.Range("K2").Value = 0.7379077026284397
$ws.Range("L2").Value = 0.7379077026284397
$ws.Range("M2").Value = 0.7400977401822874
$ws.Range("N2").Value = 0.7400977401822874
$ws.Range("O2").Value = 0.7418283274400723
$ws.Range("P2").Value = 0.8967660060965225
$ws.Range("Q2").Value = 0.8967660060965225
$ws.Range("R2").Value = 0.9473992434676685
$ws.Range("S2").Value = 0.9473992434676685
$ws.Range("T2").Value = 0.9521121941612599
$ws.Range("U2").Value = 0.9700437909678342
$ws.Range("V2").Value = 0.9700437909678342
$ws.Range("W2").Value = 0.9700437909678342
$ws.Range("X2").Value = 0.97500497814269
$ws.Range("Y2").Value = 0.97500497814269
$ws.Range("Z2").Value = 0.9959426637944316
$ws.Range("AA2").Value = 0.9959426637944316
$ws.Range("AB2").Value = 0.9959426637944316
$ws.Range("AC2").Value = 0.9959426637944316
$ws.Range("AD2").Value = 1
$ws.Range("AE2").Value = 1
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 1
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 1
$ws.Range("AJ2").Value = 1
$ws.Range("AK2").Value = 1
$ws.Range("D3").Value = 0.007761633782800323
$ws.Range("E3").Value = 0.1820224167296071
$ws.Range("F3").Value = 0.5043151789925492
$ws.Range("G3").Value = 0.5883212682088932
$ws.Range("H3").Value = 0.5927699824251023
$ws.Range("I3").Value = 0.6092951067238218
$ws.Range("J3").Value = 0.6092951067238218
$ws.Range("K3").Value = 0.6276396568762921
$ws.Range("L3").Value = 0.6281514751625253
$ws.Range("M3").Value = 0.6456452519009499
$ws.Range("N3").Value = 0.6456452519009499
$ws.Range("O3").Value = 0.6639711794755927
$ws.Range("P3").Value = 0.8317191685531655
$ws.Range("Q3").Value = 0.8317191685531655
$ws.Range("R3").Value = 0.9017294337262266
$ws.Range("S3").Value = 0.9140916523820335
$ws.Range("T3").Value = 0.9140916523820335
$ws.Range("U3").Value = 0.961352246734309
$ws.Range("V3").Value = 0.961352246734309
$ws.Range("W3").Value = 0.961352246734309
$ws.Range("X3").Value = 0.9710940092169983
$ws.Range("Y3").Value = 0.9710940092169983
$ws.Range("Z3").Value = 0.9914566738035794
$ws.Range("AA3").Value = 0.9914566738035794
$ws.Range("AB3").Value = 0.9914566738035794
$ws.Range("AC3").Value = 0.9914566738035794
$ws.Range("AD3").Value = 0.9990295627598408
$ws.Range("F4").Value = 0.1394334263860402
$ws.Range("G4").Value = 0.1394334263860402
$ws.Range("H4").Value = 0.4459735309775735
$ws.Range("I4").Value = 0.4467976507679441
$ws.Range("J4").Value = 0.4471246593851937
$ws.Range("K4").Value = 0.4703075309223458
$ws.Range("L4").Value = 0.4703075309223458
$ws.Range("M4").Value = 0.4785292087295791
$ws.Range("N4").Value = 0.4785292087295791
$ws.Range("O4").Value = 0.5200467357888819
$ws.Range("P4").Value = 0.5200467357888819
$ws.Range("Q4").Value = 0.6390062805812196
$ws.Range("R4").Value = 0.7755102523497603
$ws.Range("S4").Value = 0.8108853589486913
$ws.Range("T4").Value = 0.8270302779727297
$ws.Range("U4").Value = 0.8952599424381498
$ws.Range("V4").Value = 0.8952599424381498
$ws.Range("W4").Value = 0.926842632703224
$ws.Range("X4").Value = 0.926842632703224
$ws.Range("Y4").Value = 0.926842632703224
$ws.Range("Z4").Value = 0.9530681041768566
$ws.Range("AA4").Value = 0.9530681041768566
$ws.Range("AB4").Value = 0.9530681041768566
$ws.Range("AC4").Value = 0.9530681041768566
$ws.Range("AD4").Value = 0.9530681041768566
$ws.Range("AE4").Value = 0.9678524557409328
$ws.Range("AF4").Value = 1
$ws.Range("AG4").Value = 1
$ws.Range("AH4").Value = 1
$ws.Range("AI4").Value = 1
$ws.Range("AJ4").Value = 1
$ws.Range("AK4").Value = 1
$ws.Range("F5").Value = 0.203605394793228
$ws.Range("G5").Value = 0.203605394793228
$ws.Range("H5").Value = 0.5118630776288123
$ws.Range("I5").Value = 0.5221664754563949
$ws.Range("J5").Value = 0.5598968036807404
$ws.Range("K5").Value = 0.58945800809649
$ws.Range("L5").Value = 0.58945800809649
$ws.Range("M5").Value = 0.5904938234996756
$ws.Range("N5").Value = 0.5957004136148842
$ws.Range("O5").Value = 0.6059911032471067
$ws.Range("P5").Value = 0.6059911032471067
$ws.Range("Q5").Value = 0.6905949303301421
$ws.Range("R5").Value = 0.771695487858725
$ws.Range("S5").Value = 0.8160173538106834
$ws.Range("T5").Value = 0.8217525985828542
$ws.Range("U5").Value = 0.9200292269191265
$ws.Range("V5").Value = 0.9200292269191265
$ws.Range("W5").Value = 0.9292167517137494
$ws.Range("X5").Value = 0.9466778788398902
$ws.Range("Y5").Value = 0.9466778788398902
$ws.Range("Z5").Value = 0.958676217617338
$ws.Range("AA5").Value = 0.9643105885594532
$ws.Range("AB5").Value = 0.9759920347057951
$ws.Range("AC5").Value = 0.9759920347057951
$ws.Range("AD5").Value = 0.9759920347057951
$ws.Range("AE5").Value = 0.9759920347057951
$ws.Range("AF5").Value = 0.9999999999999999
$ws.Range("AG5").Value = 0.9999999999999999
$ws.Range("AH5").Value = 0.9999999999999999
$ws.Range("AI5").Value = 0.9999999999999999
$ws.Range("AJ5").Value = 0.9999999999999999
$ws.Range("AK5").Value = 0.9999999999999999
$ws.Range("E6").Value = 0.2744254467974087
$ws.Range("F6").Value = 0.3124871335555507
$ws.Range("G6").Value = 0.43505250727291
$ws.Range("H6").Value = 0.43505250727291
$ws.Range("I6").Value = 0.4691638769229716
$ws.Range("J6").Value = 0.4984378379214139
$ws.Range("K6").Value = 0.5173539743995128
$ws.Range("L6").Value = 0.5173539743995128
$ws.Range("M6").Value = 0.5173539743995128
$ws.Range("N6").Value = 0.5284141675363246
$ws.Range("O6").Value = 0.5284141675363246
$ws.Range("P6").Value = 0.6667712348184762
$ws.Range("Q6").Value = 0.7006858056208057
$ws.Range("R6").Value = 0.8280493886173371
$ws.Range("S6").Value = 0.8280493886173371
$ws.Range("T6").Value = 0.9659294144108383
$ws.Range("U6").Value = 0.9891888572093006
$ws.Range("V6").Value = 0.9956840754679194
$ws.Range("W6").Value = 0.9963596230316747
$ws.Range("X6").Value = 0.9963596230316747
$ws.Range("Y6").Value = 0.9963596230316747
$ws.Range("Z6").Value = 0.9963596230316747
$ws.Range("AA6").Value = 0.9963596230316747
$ws.Range("AB6").Value = 0.9963596230316747
$ws.Range("AC6").Value = 0.9963596230316747
$ws.Range("AD6").Value = 0.9963596230316747

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.5892774124955195
$ws.Range("F3").Value = 0.5043151789925492
$ws.Range("D4").Value = 14
$ws.Range("F4").Value = 0.5200467357888819
$ws.Range("G4").Value = 11
$ws.Range("F5").Value = 0.5118630776288123
$ws.Range("D6").Value = 10
$ws.Range("F6").Value = 0.5173539743995128
$ws.Range("G6").Value = 8

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F2").Value = 0.7038137204085497
$ws.Range("D3").Value = 15
$ws.Range("F3").Value = 0.8317191685531655
$ws.Range("G3").Value = 13
$ws.Range("D4").Value = 17
$ws.Range("F4").Value = 0.7755102523497603
$ws.Range("G4").Value = 14
$ws.Range("D5").Value = 17
$ws.Range("F5").Value = 0.771695487858725
$ws.Range("G5").Value = 14
$ws.Range("D6").Value = 16
$ws.Range("F6").Value = 0.7006858056208057
$ws.Range("G6").Value = 14

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 15
$ws.Range("F2").Value = 0.8967660060965225
$ws.Range("G2").Value = 13
$ws.Range("F3").Value = 0.8317191685531655
$ws.Range("D4").Value = 18
$ws.Range("F4").Value = 0.8108853589486913
$ws.Range("G4").Value = 15
$ws.Range("D5").Value = 18
$ws.Range("F5").Value = 0.8160173538106834
$ws.Range("G5").Value = 15
$ws.Range("F6").Value = 0.8280493886173371

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 17
$ws.Range("F2").Value = 0.9473992434676685
$ws.Range("G2").Value = 15
$ws.Range("D3").Value = 17
$ws.Range("F3").Value = 0.9017294337262266
$ws.Range("G3").Value = 15
$ws.Range("D4").Value = 22
$ws.Range("F4").Value = 0.926842632703224
$ws.Range("G4").Value = 19
$ws.Range("F5").Value = 0.9200292269191265
$ws.Range("F6").Value = 0.9659294144108383
